# Insert a new client row ("BRITO CARDENAS RUTH CECILIA") as the 3rd
# client (row 4) on both sheets, pushing the existing rows (and the
# trailing summary row) down by one, and refresh the "0 de N" counter
# text on the first sheet's summary row to reflect the new total.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" --------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(4).Insert()

$ws1.Range("A4").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B4").Value = "BRITO CARDENAS RUTH CECILIA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(4, $col).Value = 0
}

# Summary row (previously row 8, now row 9) - bump the "0 de 6" counter
# to "0 de 7" now that there are 7 client rows.
$ws1.Range("C9:R9").Value = "0 de 7"

# --- Sheet 2: "VENTA MENSUAL" ------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(4).Insert()

$ws2.Range("A4").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B4").Value = "BRITO CARDENAS RUTH CECILIA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(4, $col).Value = 0
}
